# Update cryptos list values to match the latest scrape (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.510.82'
$ws.Range("D3").Value = '3.672.61'
$ws.Range("E3").Value = '  -0.94%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '622.04'
$ws.Range("E5").Value = '  -7.62%  '
$ws.Range("D6").Value = '158.90'
$ws.Range("E6").Value = '  -1.94%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  -0.48%  '
$ws.Range("E9").Value = '  -1.50%  '
$ws.Range("E10").Value = '  +1.50%  '
$ws.Range("E11").Value = '  -0.89%  '
$ws.Range("E12").Value = '  -2.82%  '
$ws.Range("D13").Value = '4.293.13'
$ws.Range("E13").Value = '  -0.89%  '
$ws.Range("D14").Value = '32.24'
$ws.Range("E14").Value = '  -1.99%  '
$ws.Range("D15").Value = '3.675.14'
$ws.Range("E15").Value = '  -1.11%  '
$ws.Range("D16").Value = '69.514.76'
$ws.Range("E16").Value = '  -0.24%  '
$ws.Range("E17").Value = '  +0.58%  '
$ws.Range("D18").Value = '15.89'
$ws.Range("E18").Value = '  -2.56%  '
$ws.Range("D19").Value = '6.49'
$ws.Range("E19").Value = '  -0.25%  '
$ws.Range("D20").Value = '10.26'
$ws.Range("E20").Value = '  +4.52%  '
$ws.Range("D21").Value = '468.55'
$ws.Range("E21").Value = '  -1.28%  '
$ws.Range("D22").Value = '0.650'
$ws.Range("E22").Value = '  -0.76%  '
$ws.Range("D23").Value = '79.59'
$ws.Range("E23").Value = '  -1.03%  '
$ws.Range("D24").Value = '3.820.21'
$ws.Range("E24").Value = '  -0.90%  '
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("D26").Value = '11.15'
$ws.Range("E26").Value = '  +1.16%  '
$ws.Range("D27").Value = '0.0000121'
$ws.Range("E27").Value = '  -5.56%  '
$ws.Range("D28").Value = '8.62'
$ws.Range("E28").Value = '  -5.97%  '
$ws.Range("D29").Value = '2.61'
$ws.Range("E29").Value = '  -3.25%  '
$ws.Range("D30").Value = '1.65'
$ws.Range("E30").Value = '  -4.23%  '
$ws.Range("E31").Value = '  +0.25%  '
$ws.Range("D32").Value = '1.96'
$ws.Range("E32").Value = '  -2.60%  '
$ws.Range("D33").Value = '26.60'
$ws.Range("E33").Value = '  -1.36%  '
$ws.Range("B34").Value = 'RenzoRestakedETH'
$ws.Range("C34").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D34").Value = '3.683.00'
$ws.Range("E34").Value = '  -0.39%  '
$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").Value = '6.38'
$ws.Range("E35").Value = '  -3.28%  '
$ws.Range("E36").Value = '  -4.11%  '
$ws.Range("D37").Value = '8.26'
$ws.Range("E37").Value = '  -3.29%  '
$ws.Range("E38").Value = '  +0.03%  '
$ws.Range("D39").Value = '178.32'
$ws.Range("E39").Value = '  +2.55%  '
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("E41").Value = '  -2.27%  '
$ws.Range("D42").Value = '5.78'
$ws.Range("E42").Value = '  -5.64%  '
$ws.Range("D43").Value = '0.0892'
$ws.Range("E43").Value = '  -2.62%  '
$ws.Range("D44").Value = '0.923'
$ws.Range("E44").Value = '  -1.99%  '
$ws.Range("D45").Value = '29.28'
$ws.Range("E45").Value = '  +5.82%  '
$ws.Range("D46").Value = '46.75'
$ws.Range("E46").Value = '  -0.73%  '
$ws.Range("D47").Value = '2.69'
$ws.Range("E47").Value = '  -2.84%  '
$ws.Range("D48").Value = '7.84'
$ws.Range("E48").Value = '  -0.50%  '
$ws.Range("D49").Value = '0.000263'
$ws.Range("E49").Value = '  -6.23%  '
$ws.Range("E50").Value = '  -5.64%  '
$ws.Range("E51").Value = '  -3.21%  '
